$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.598.42"
$ws.Range("E2").Value = "  +1.49%  "

# Row 3
$ws.Range("D3").Value = "1.826.60"
$ws.Range("E3").Value = "  +1.22%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.13%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5311"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3978"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.70%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07749"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.60%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.07%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.118"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.85%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.68%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.319"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.51%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.578"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.82%  "

# Row 15
$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.002"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.15%  "

# Row 16
$ws.Range("D16").Value = "1.828.63"
$ws.Range("E16").Value = "  +1.32%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.42%  "

# Row 18
$ws.Range("E18").Value = "  +2.08%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06615"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.23%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.70%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.089"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.51%  "

# Row 23
$ws.Range("D23").Value = "28.610.43"
$ws.Range("E23").Value = "  +1.44%  "

# Row 24
$ws.Range("E24").Value = "  -0.18%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.234"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.62%  "

# Row 26
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.53%  "

# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.37%  "

# Row 28
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "2.040.48"
$ws.Range("E28").Value = "  +1.40%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.416"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.90%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.70%  "

# Row 31
$ws.Range("E31").Value = "  +2.62%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1124"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.08%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.738"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.86%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.659"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.07%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07322"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.16%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2270"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.73%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02352"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.72%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.910"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.83%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.199"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.22%  "

# Row 40
$ws.Range("E40").Value = "  +1.73%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6296"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.58%  "

# Row 42
$ws.Range("E42").Value = "  +1.78%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.09%  "

# Row 44
$ws.Range("E44").Value = "  -1.49%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.97%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5938"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.90%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.722"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.13%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.49%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.66%  "

# Row 50
$ws.Range("E50").Value = "  +0.11%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06954"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.89%  "
